$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New staff row: "cake" / "33" appended at row 20 (matching existing A16:B17 layout)
$ws.Range("A20").Value = "cake"

# Writing "33" via .Value would be auto-detected as a number, and forcing
# text via NumberFormat/quote-prefix would introduce a brand new cell
# style. Instead, stage the text value as a formula result (style-neutral)
# in a scratch cell, then copy only the *value* (PasteSpecial values) into
# B20 so it lands as a shared-string text cell using the default style.
$scratch = $ws.Range("D1")
$scratch.Formula = "=""33"""
$scratch.Copy()
$ws.Range("B20").PasteSpecial(-4163)
$scratch.Clear()
